# IPDO.xlsx - add a new "22 Maio de 2016" data row (row 63) to the
# "Tabela1" sheet, following the exact same layout used by the rows
# above it (one row per report date, columns A:I = "Programado",
# columns J:Q = "Verificado").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabela1")

$row = 63

# --- Programado -----------------------------------------------------
$ws.Cells.Item($row, 1).Value = " 22 Maio de 2016"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "29.037"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "9.653"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "1.990"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 5).Value = "8.234"
$ws.Cells.Item($row, 5).Style = "Normal"

$ws.Cells.Item($row, 6).NumberFormat = "@"
$ws.Cells.Item($row, 6).Value = "3.556"
$ws.Cells.Item($row, 6).Style = "Normal"

$ws.Cells.Item($row, 7).NumberFormat = "@"
$ws.Cells.Item($row, 7).Value = "52.470"
$ws.Cells.Item($row, 7).Style = "Normal"

$ws.Cells.Item($row, 8).NumberFormat = "@"
$ws.Cells.Item($row, 8).Value = "0"
$ws.Cells.Item($row, 8).Style = "Normal"

$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 9).Value = "52.470"
$ws.Cells.Item($row, 9).Style = "Normal"

# --- Verificado -------------------------------------------------------
$ws.Cells.Item($row, 10).NumberFormat = "@"
$ws.Cells.Item($row, 10).Value = "29.472"
$ws.Cells.Item($row, 10).Style = "Normal"

$ws.Cells.Item($row, 11).NumberFormat = "@"
$ws.Cells.Item($row, 11).Value = "9.636"
$ws.Cells.Item($row, 11).Style = "Normal"

$ws.Cells.Item($row, 12).NumberFormat = "@"
$ws.Cells.Item($row, 12).Value = "2.015"
$ws.Cells.Item($row, 12).Style = "Normal"

$ws.Cells.Item($row, 13).NumberFormat = "@"
$ws.Cells.Item($row, 13).Value = "7.888"
$ws.Cells.Item($row, 13).Style = "Normal"

$ws.Cells.Item($row, 14).NumberFormat = "@"
$ws.Cells.Item($row, 14).Value = "3.313"
$ws.Cells.Item($row, 14).Style = "Normal"

$ws.Cells.Item($row, 15).NumberFormat = "@"
$ws.Cells.Item($row, 15).Value = "52.324"
$ws.Cells.Item($row, 15).Style = "Normal"

$ws.Cells.Item($row, 16).NumberFormat = "@"
$ws.Cells.Item($row, 16).Value = "0"
$ws.Cells.Item($row, 16).Style = "Normal"

$ws.Cells.Item($row, 17).NumberFormat = "@"
$ws.Cells.Item($row, 17).Value = "52.324"
$ws.Cells.Item($row, 17).Style = "Normal"

# Column R (18) never carried real data for any row in this table - it
# holds the same "blank" placeholder cell that used to sit at the end
# of the previous last row (R62). Move that placeholder down to the
# newly appended row (R63) so the table keeps exactly one trailing
# blank marker, just like before the row was added.
$placeholder = $ws.Cells.Item(62, 18).Value()
$ws.Cells.Item(62, 18).ClearContents()
$ws.Cells.Item($row, 18).NumberFormat = "@"
$ws.Cells.Item($row, 18).Value = $placeholder
$ws.Cells.Item($row, 18).Style = "Normal"
